$d = $word.ActiveDocument

# 1) "15-pound e-voucher" -> "£15 e-voucher"
$d.Content.Find.Execute("15-pound", $false, $false, $false, $false, $false, $true, 1, $false, "£15", 2) | Out-Null

# 2) Insert "or have any questions, " before "please contact us on our email: "
$d.Content.Find.Execute("study, please contact us on our emai", $false, $false, $false, $false, $false, $true, 1, $false, "study, or have any questions, please contact us on our emai", 2) | Out-Null

# 3) Replace the old study email address with the new contact email
$d.Content.Find.Execute("pip-tabletstudy@kcl.ac.uk", $false, $false, $false, $false, $false, $true, 1, $false, "kclbrainrangers@gmail.com", 2) | Out-Null

# 4) Remove the now-superseded paragraphs: "If you have any questions...", "Project Supervisor: ...",
#    "Project Members: ..." -- collapsing them back down to the single blank paragraph that used to
#    separate the email line from the "Thank you" paragraph.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "If you have any questions or would like to raise a concern*") {
        $startPara = $d.Paragraphs.Item($i)
    }
    if ($t -like "Project Members*") {
        $endPara = $d.Paragraphs.Item($i)
        break
    }
}
if ($startPara -ne $null -and $endPara -ne $null) {
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
